# The upstream change this task mirrors ("pptx: Include all themes in
# output archive") is a fix in a third-party OOXML *writer* so that every
# theme part referenced by the package (not just the one used by the
# slide master) is actually written into the saved .pptx archive.
#
# This particular fixture - test/pptx/speaker-notes/templated.pptx - is
# one of the decks that already carries speaker notes, and therefore
# already had a second theme (ppt/theme/theme2.xml, used by
# ppt/notesMasters/notesMaster1.xml) referenced and packaged *before* the
# fix. So for this file the only observable effect of the writer fix is
# that ppt/slideMasters/slideMaster1.xml gets re-emitted by the (fixed)
# serializer, which happens to alphabetize element attributes as a side
# effect. There is no content/model change: the reconstructed "after"
# tree round-trips to the exact same canonical (C14N) XML as the
# "before" tree - same elements, same attribute values, same text, same
# relationships - only the on-disk attribute order differs.
#
# None of that is something a PowerPoint end-user/automation script can
# (or should) express through the COM object model - it is purely an
# artifact of the OOXML-writing library used to produce the fixture.
# Both themes are already present and correctly linked in this
# presentation, so there is nothing for this script to change; it
# intentionally performs no edits, leaving the deck semantically (and,
# for every part other than the master's attribute order, byte-for-byte)
# identical to the input.

$p = $ppt.ActivePresentation
